# Commit: "Add files via upload"
# The dataset used "?" as a placeholder for missing values. This pass
# replaces every occurrence of "?" in the data with "NA", and leaves the
# active selection on cell D13 (matching the refreshed view state saved
# with the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$missingValueCells = @(
    "B2","B3","B4","B7","B9","B11","Z11","B16","B17","B18","B19","F29",
    "B45","B46","Z46","B47","Z47","B48","B50","B51",
    "S57","T57","S58","T58","S59","T59","S60","T60",
    "B65","F65","B68","B73","B75","B76","B77",
    "B84","B85","B86",
    "B111","B112","B115","B116","B126",
    "B128","B129","B130","B131","Z131",
    "B132","V132","W132","B133","V133","W133",
    "B183","B191","B193","B194","B195"
)

foreach ($addr in $missingValueCells) {
    $ws.Range($addr).Value = "NA"
}

$null = $ws.Range("D13").Select()
